$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5942.2856
$ws.Range("I100").Value = 1932.6666
$ws.Range("J100").Value = 30000
$ws.Range("K100").Value = 1932.6666
$ws.Range("L100").Value = 30000
$ws.Range("M100").Value = -1391.6666
$ws.Range("N100").Value = -31082
$ws.Range("H103").Value = 995.9
$ws.Range("J103").Value = 996
$ws.Range("L103").Value = 2988
$ws.Range("N103").Value = -4160
$ws.Range("H135").Value = 827.7727
$ws.Range("I135").Value = 771.9524
$ws.Range("K135").Value = 6947.5716
$ws.Range("M135").Value = -4412.5716
$ws.Range("H140").Value = 125000
$ws.Range("J140").Value = 125000
$ws.Range("L140").Value = 125000
$ws.Range("N140").Value = -135360
$ws.Range("H141").Value = 2833.4583
$ws.Range("I141").Value = 2428.8096
$ws.Range("K141").Value = 7286.4288
$ws.Range("M141").Value = -2106.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1520.125
$ws.Range("I61").Value = 1520.125
$ws.Range("K61").Value = 1520.125
$ws.Range("M61").Value = -1308.125
$ws.Range("H97").Value = 1472
$ws.Range("J97").Value = 3592.3333
$ws.Range("L97").Value = 3592.3333
$ws.Range("N97").Value = -4584.3333
$ws.Range("H122").Value = 2458.6924
$ws.Range("I122").Value = 2458.6924
$ws.Range("K122").Value = 7376.0772
$ws.Range("M122").Value = -4926.0772
$ws.Range("H132").Value = 2743.8147
$ws.Range("I132").Value = 1920.7222
$ws.Range("K132").Value = 5762.1666
$ws.Range("M132").Value = -3232.1666
$ws.Range("H136").Value = 1520.125
$ws.Range("I136").Value = 1520.125
$ws.Range("K136").Value = 4560.375
$ws.Range("M136").Value = -2010.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 880.9231
$ws.Range("J80").Value = 794.375
$ws.Range("L80").Value = 794.375
$ws.Range("N80").Value = -2790.375
$ws.Range("H83").Value = 880.9231
$ws.Range("J83").Value = 794.375
$ws.Range("L83").Value = 3971.875
$ws.Range("N83").Value = -13955.875
$ws.Range("H131").Value = 24999.5
$ws.Range("J131").Value = 24999.5
$ws.Range("L131").Value = 24999.5
$ws.Range("N131").Value = -35079.5
$ws.Range("H134").Value = 3138.375
$ws.Range("I134").Value = 3138.375
$ws.Range("K134").Value = 9415.125
$ws.Range("M134").Value = -6880.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3498.7144
$ws.Range("I31").Value = 3298.6
$ws.Range("K31").Value = 3298.6
$ws.Range("M31").Value = -3003.6
$ws.Range("H34").Value = 3498.7144
$ws.Range("I34").Value = 3298.6
$ws.Range("K34").Value = 3298.6
$ws.Range("M34").Value = -3096.6
$ws.Range("H58").Value = 2485
$ws.Range("I58").Value = 2132.3333
$ws.Range("J58").Value = 3014
$ws.Range("K58").Value = 2132.3333
$ws.Range("L58").Value = 3014
$ws.Range("M58").Value = -1929.3333
$ws.Range("N58").Value = -3420
$ws.Range("H86").Value = 35369.668
$ws.Range("J86").Value = 47720.332
$ws.Range("L86").Value = 47720.332
$ws.Range("N86").Value = -49966.332
$ws.Range("H89").Value = 35369.668
$ws.Range("J89").Value = 47720.332
$ws.Range("L89").Value = 238601.66
$ws.Range("N89").Value = -249833.66
$ws.Range("H132").Value = 3032.45
$ws.Range("I132").Value = 2816.0667
$ws.Range("J132").Value = 3681.6
$ws.Range("K132").Value = 8448.2001
$ws.Range("L132").Value = 11044.8
$ws.Range("M132").Value = -5918.2001
$ws.Range("N132").Value = -16104.8
$ws.Range("H134").Value = 3311.2856
$ws.Range("I134").Value = 3311.2856
$ws.Range("K134").Value = 9933.856800000001
$ws.Range("M134").Value = -7398.856800000001
$ws.Range("H136").Value = 2485
$ws.Range("I136").Value = 2132.3333
$ws.Range("J136").Value = 3014
$ws.Range("K136").Value = 6396.999899999999
$ws.Range("L136").Value = 9042
$ws.Range("M136").Value = -3846.999899999999
$ws.Range("N136").Value = -14142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 252721.75
$ws.Range("I80").Value = 3629.3333
$ws.Range("K80").Value = 10887.9999
$ws.Range("M80").Value = -9951.999899999999
$ws.Range("H83").Value = 252721.75
$ws.Range("I83").Value = 3629.3333
$ws.Range("K83").Value = 32663.9997
$ws.Range("M83").Value = -27983.9997
$ws.Range("H131").Value = 1516.6
$ws.Range("I131").Value = 955.625
$ws.Range("J131").Value = 1780.5883
$ws.Range("K131").Value = 2866.875
$ws.Range("L131").Value = 5341.7649
$ws.Range("M131").Value = 2173.125
$ws.Range("N131").Value = -15421.7649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6819.8887
$ws.Range("I80").Value = 5229.1665
$ws.Range("J80").Value = 10001.333
$ws.Range("K80").Value = 5229.1665
$ws.Range("L80").Value = 10001.333
$ws.Range("M80").Value = -4231.1665
$ws.Range("N80").Value = -11997.333
$ws.Range("H83").Value = 6819.8887
$ws.Range("I83").Value = 5229.1665
$ws.Range("J83").Value = 10001.333
$ws.Range("K83").Value = 26145.8325
$ws.Range("L83").Value = 50006.665
$ws.Range("M83").Value = -21153.8325
$ws.Range("N83").Value = -59990.665
$ws.Range("H107").Value = 540.1429000000001
$ws.Range("I107").Value = 520.2
$ws.Range("K107").Value = 520.2
$ws.Range("M107").Value = 1399.8
$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3025.5
$ws.Range("I7").Value = 3025.5
$ws.Range("K7").Value = 3025.5
$ws.Range("M7").Value = -2913.5
$ws.Range("H126").Value = 3025.5
$ws.Range("I126").Value = 3025.5
$ws.Range("K126").Value = 9076.5
$ws.Range("M126").Value = -6606.5
$ws.Range("H132").Value = 3343.3462
$ws.Range("I132").Value = 2495.1538
$ws.Range("K132").Value = 7485.4614
$ws.Range("M132").Value = -4955.4614
$ws.Range("H136").Value = 3312.5217
$ws.Range("I136").Value = 3272.2632
$ws.Range("K136").Value = 9816.7896
$ws.Range("M136").Value = -7266.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 492.83334
$ws.Range("I107").Value = 431.8
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 1295.4
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = 624.5999999999999
$ws.Range("N107").Value = -6234
$ws.Range("H113").Value = 802.63635
$ws.Range("I113").Value = 652.9
$ws.Range("K113").Value = 1958.7
$ws.Range("M113").Value = 211.3000000000002
$ws.Range("H122").Value = 1399.3334
$ws.Range("I122").Value = 1399.3334
$ws.Range("K122").Value = 4198.0002
$ws.Range("M122").Value = -1748.0002
$ws.Range("H136").Value = 1156.3226
$ws.Range("I136").Value = 851.6786
$ws.Range("K136").Value = 2555.0358
$ws.Range("M136").Value = -5.035799999999654
